# Project: LinkedIn Visuals  Script: findEmployeeCount
#
# Refresh the scraped LinkedIn numbers on the "Data" sheet:
#   - EmployeeRange (col D) "10001" -> "10001+" (and the TresVista 1,001-5,000
#     row is left alone)
#   - EmployeesOnLinkedIn (col E) gets the latest headcount snapshot
#   - the stray trailing blank row is removed
#   - the saved selection is moved off the old (now out-of-range) cell

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# The new EmployeesOnLinkedIn figures look like plain numbers, but the sheet
# stores this column as text (matching the existing "1,001-5,000" style
# entries next to it), so force Text format before typing them in -- otherwise
# they'd land as numeric literals instead of shared strings.
$ws.Range("D3:F6").NumberFormat = "@"

$ws.Range("D3").Value() = "10001+"
$ws.Range("E3").Value() = "15515"

$ws.Range("D4").Value() = "1,001-5,000"
$ws.Range("E4").Value() = "139"

$ws.Range("D5").Value() = "10001+"
$ws.Range("E5").Value() = "57501"

$ws.Range("D6").Value() = "10001+"
$ws.Range("E6").Value() = "87651"

# Re-apply the plain bordered look used by the rest of the table (this also
# clears the wrapText formatting the EmployeeRange/EmployeesOnLinkedIn/
# LinkedinLink columns had before).
$ws.Range("C3").Copy()
$ws.Range("D3:F6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Drop the now-unused trailing blank row.
$ws.Range("A7:B7").EntireRow.Delete()

# Move the active cell selection off the stale reference.
$ws.Range("E12").Select()
